# Realestate Update resale numbers 2024-01-23 14:25
# Appends the new daily resale-number observation as row 87 on the
# (single) active "CityResaleNum" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

# Columns A (Date) and D (Week) hold text that looks numeric/date-like
# ("2024-01-23", "03"). Pre-format them as Text so Excel's COM layer
# stores them verbatim instead of silently coercing them into a date
# serial / integer, matching how every other row in the column is
# stored.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("D$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2024-01-23"
$ws.Range("B$row").Value = "14:25:07"
$ws.Range("C$row").Value = "Tuesday"
$ws.Range("D$row").Value = "03"
$ws.Range("E$row").Value = 138481
$ws.Range("F$row").Value = 141269
$ws.Range("G$row").Value = 171036
$ws.Range("H$row").Value = 148800
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 123104
$ws.Range("K$row").Value = 223626
$ws.Range("L$row").Value = 255955
$ws.Range("M$row").Value = 185250
$ws.Range("N$row").Value = 110238
$ws.Range("O$row").Value = 41289
$ws.Range("P$row").Value = 30877
$ws.Range("Q$row").Value = 73448
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 42326
$ws.Range("T$row").Value = -1
